# Warheads.xlsx - buff/debuff + hit-multiplicity + hitpoint scale unification
# See commit message:
#   1. Implemented hitpoint buffs/debuffs for ship hull, components, turrets.
#   2. Implemented hit multiplicity (n rolls of armour penetration -> n damage).
#   3. Unified scale of ship hitpoints with strike craft hitpoints (x10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Base weapon rows (J = HD column, raw values; T = HD/Interval derived) ---
# Light Autocannon
$ws.Range("J2").Value = 20
$ws.Range("J3").Value = 40
$ws.Range("J4").Value = 10

# Light HVGun
$ws.Range("J6").Value = 80
$ws.Range("J7").Value = 160
$ws.Range("J8").Value = 60

# Light Howitzer
$ws.Range("J10").Value = 60
$ws.Range("J11").Value = 240
$ws.Range("J12").Value = 60

# Light strike-craft weapon rows (plain values, no formulas)
$ws.Range("J15").Value = 80
$ws.Range("J16").Value = 100
$ws.Range("J17").Value = 120

# --- Medium Autocannon (row 20-22): J = 2*T{row-18}*D{row}; row21 multiplier bumped 2 -> 2.5
$ws.Range("J21").Formula = "=2.5*T3*D21"

# --- Medium HVGun (row 24-26): I24 multiplier bumped 2 -> 3; J25 multiplier bumped 2 -> 2.5
$ws.Range("I24").Formula = "=3*S6*D24"
$ws.Range("J25").Formula = "=2.5*T7*D25"

# --- Medium Howitzer (row 28-30): J29 multiplier bumped 2 -> 2.5
$ws.Range("J29").Formula = "=2.5*T11*D29"

# --- Heavy Autocannon (row 38-40): I38 multiplier bumped 2.5 -> 3; J39 multiplier bumped 2.5 -> 3
$ws.Range("I38").Formula = "=3*S20*D38"
$ws.Range("J39").Formula = "=3*T21*D39"

# --- Heavy HVGun (row 42-44): I42 multiplier bumped 2.5 -> 3.5; J43 multiplier bumped 2.5 -> 3
$ws.Range("I42").Formula = "=3.5*S24*D42"
$ws.Range("J43").Formula = "=3*T25*D43"

# --- Heavy Howitzer (row 46-48): I46 multiplier bumped 2.5 -> 3; J47 multiplier bumped 2.5 -> 3.5
$ws.Range("I46").Formula = "=3*S28*D46"
$ws.Range("J47").Formula = "=3.5*T29*D47"

# --- StrikeCraft row 56: J56 loses its formula, becomes a flat literal
$ws.Range("J56").Value = 50

# --- StrikeCraft row 57
$ws.Range("J57").Value = 2

# --- Torpedo rows 60-62: damage-share multipliers rebalanced
$ws.Range("H60").Formula = "=R29*D60*0.25"
$ws.Range("I60").Formula = "=S29*D60*0.65"
$ws.Range("J60").Formula = "=T29*D60*0.5"

$ws.Range("H61").Formula = "=R47*D61*0.25"
$ws.Range("I61").Formula = "=S47*D61*0.3"
$ws.Range("J61").Formula = "=T47*D61*0.25"

$ws.Range("H62").Formula = "=R29*D60*0.25"
$ws.Range("I62").Formula = "=S29*D60*0.65"
$ws.Range("J62").Formula = "=T29*D60*0.45"

# --- View state: scroll so row 49 is at top, select J58 (matches author's saved cursor position)
try { $excel.ActiveWindow.ScrollRow = 49 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 1 } catch {}
$ws.Range("J58").Select()
